$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space character used in the original site copy.
$nbsp = [char]0x00A0

# 1) Delete row 80 ("li: St John's School, MITCHAM") which was removed outright.
#    This shifts every row below it up by one, matching rows 80-247 of the target.
$ws.Rows.Item(80).Delete() | Out-Null

# 2) Update the "On this page" banner (row 7) with the new day/date/time.
$row7 = "On this page" + `
    "Current school and early childhood service, TAFE closures and relocations:" + `
    "Bus service cancellations or alterations" + `
    "Current school" + $nbsp + "and early childhood service," + $nbsp + `
    "TAFE closures and relocations for Thursday 1" + $nbsp + "October," + $nbsp + `
    "(as at 10:00am, 1" + $nbsp + "October)" + `
    "South-Eastern Victoria RegionEarly childhood services"
$ws.Range("A7").Value = $row7

# 3) Row 78: the closed-school list entry was removed, leaving just the lead-in text.
$ws.Range("A78").Value = "Schools closedThe Department hasbeen advised of the following school closures:"

# 4) Row 79: the single remaining closed school changed.
$ws.Range("A79").Value = "li: St John's School, FRANKSTON EAST"

# 5) Row 242 (formerly 243, now shifted up after the deletion): updated "Last Update" stamp.
$ws.Range("A242").Value = "li: Last Update: 01 October 2020"
